# redid script and added comments
# Update "Number of Mitophagy Events Outside Dialated Area" (column E) values
# for several rows based on re-analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7   = 5
    17  = 5
    40  = 8
    47  = 4
    54  = 5
    65  = 6
    82  = 1
    86  = 4
    99  = 1
    103 = 4
    114 = 4
    122 = 1
    130 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
